$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 date-like text value ("08/08/2025") must stay a plain text cell
# (matching the rest of the sheet, which stores dates as inline strings,
# not real Excel dates). Forcing the NumberFormat to Text ("@") before
# assigning the value prevents Excel's autodetection from turning it into
# a date serial; ClearFormats() afterwards drops the now-unneeded explicit
# style so the cell ends up with the default (no "s" attribute), exactly
# like the sibling rows.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "08/08/2025"
$ws.Range("A21").ClearFormats()

$ws.Range("B21").Value = "Estudiantes"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "Ind. Rivadavia"
$ws.Range("F21").Value = "W"
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 1.32
$ws.Range("L21").Value = 0.97
$ws.Range("M21").Value = 11
$ws.Range("N21").Value = 7
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 4
